$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.437.17"
$ws.Range("E2").Value = "  +0.09%  "
$ws.Range("D3").Value = "3.683.52"
$ws.Range("E3").Value = "  -0.48%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "685.72"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.19%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "159.70"
$ws.Range("D6").Style = "Normal"
$ws.Range("E7").Value = "  +0.05%  "
$ws.Range("E8").Value = "  -0.19%  "
$ws.Range("E9").Value = "  -1.62%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.06"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.87%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.435"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -3.68%  "
$ws.Range("E12").Value = "  -1.86%  "
$ws.Range("D13").Value = "4.306.33"
$ws.Range("E13").Value = "  -0.29%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "32.39"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -3.84%  "
$ws.Range("D15").Value = "3.671.26"
$ws.Range("E15").Value = "  -0.21%  "
$ws.Range("D16").Value = "69.392.65"
$ws.Range("E16").Value = "  +0.04%  "
$ws.Range("E17").Value = "  +1.78%  "
$ws.Range("E18").Value = "  -3.58%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.39"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -3.48%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "471.35"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.90%  "
$ws.Range("E21").Value = "  +1.01%  "
$ws.Range("E22").Value = "  -2.21%  "
$ws.Range("E23").Value = "  +0.17%  "
$ws.Range("D24").Value = "3.832.63"
$ws.Range("E24").Value = "  -0.18%  "
$ws.Range("E26").Value = "  -3.89%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.98"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -5.62%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.20"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -3.52%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.71"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.49%  "
$ws.Range("E30").Value = "  -4.57%  "
$ws.Range("E31").Value = "  -5.62%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.56"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.73%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.999"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.03%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "26.89"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.17%  "
$ws.Range("D35").Value = "3.659.08"
$ws.Range("E35").Value = "  -0.16%  "
$ws.Range("E36").Value = "  -2.97%  "
$ws.Range("E37").Value = "  -4.38%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.12"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.27%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.22"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.62%  "
$ws.Range("E41").Value = "  -4.29%  "
$ws.Range("E42").Value = "  +0.09%  "
$ws.Range("E43").Value = "  -2.05%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "165.84"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +3.29%  "
$ws.Range("E45").Value = "  -1.53%  "
$ws.Range("B46").Value = "dogwifhat"
$ws.Range("C46").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.74"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -3.87%  "
$ws.Range("B47").Value = "FLOKI"
$ws.Range("C47").Value = "https://coinranking.com/coin/fmHk13Rqw+floki-floki"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.000281"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.19%  "
$ws.Range("E48").Value = "  +4.89%  "
$ws.Range("E49").Value = "  +0.05%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "27.62"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.18%  "
$ws.Range("E51").Value = "  -3.68%  "
